# Scheduled runner update: refresh Leve profit-calculation inputs/outputs
# (currentAveragePrice* / LevePrice* / LeveProfit* columns H:N) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 581
$ws.Range("J2").Value = 1764
$ws.Range("L2").Value = 1764
$ws.Range("N2").Value = -1990

$ws.Range("H5").Value = 505
$ws.Range("I5").Value = 505
$ws.Range("K5").Value = 505
$ws.Range("M5").Value = -390

$ws.Range("H28").Value = 652.2
$ws.Range("I28").Value = 475.4
$ws.Range("J28").Value = 1182.6
$ws.Range("K28").Value = 475.4
$ws.Range("L28").Value = 1182.6
$ws.Range("M28").Value = 9.600000000000023
$ws.Range("N28").Value = -2152.6

$ws.Range("H64").Value = 14709633
$ws.Range("I64").Value = 27781622
$ws.Range("K64").Value = 27781622
$ws.Range("M64").Value = -27781374

$ws.Range("H67").Value = 14709633
$ws.Range("I67").Value = 27781622
$ws.Range("K67").Value = 27781622
$ws.Range("M67").Value = -27780764

$ws.Range("H69").Value = 18248.25
$ws.Range("J69").Value = 19687.25
$ws.Range("L69").Value = 59061.75
$ws.Range("N69").Value = -60809.75

$ws.Range("H72").Value = 18248.25
$ws.Range("J72").Value = 19687.25
$ws.Range("L72").Value = 177185.25
$ws.Range("N72").Value = -185921.25

$ws.Range("H82").Value = 4011.6365
$ws.Range("I82").Value = 2912.8
$ws.Range("K82").Value = 8738.400000000001
$ws.Range("M82").Value = -8332.400000000001

$ws.Range("H85").Value = 4011.6365
$ws.Range("I85").Value = 2912.8
$ws.Range("K85").Value = 8738.400000000001
$ws.Range("M85").Value = -7334.400000000001

$ws.Range("H96").Value = 847.7143
$ws.Range("J96").Value = 289.75
$ws.Range("L96").Value = 869.25
$ws.Range("N96").Value = -3615.25

$ws.Range("H116").Value = 5659.9375
$ws.Range("I116").Value = 5563.8887
$ws.Range("J116").Value = 5783.4287
$ws.Range("K116").Value = 5563.8887
$ws.Range("L116").Value = 5783.4287
$ws.Range("M116").Value = -2121.8887
$ws.Range("N116").Value = -12667.4287

$ws.Range("H132").Value = 3084.625
$ws.Range("I132").Value = 3143.6
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 9430.799999999999
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -6900.799999999999
$ws.Range("N132").Value = -11660

$ws.Range("H141").Value = 1358.2858
$ws.Range("I141").Value = 1358.2858
$ws.Range("K141").Value = 4074.8574
$ws.Range("M141").Value = 1105.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5671.0713
$ws.Range("I45").Value = 8185.25
$ws.Range("K45").Value = 8185.25
$ws.Range("M45").Value = -7808.25

$ws.Range("H74").Value = 31252886
$ws.Range("I74").Value = 31252886
$ws.Range("K74").Value = 31252886
$ws.Range("M74").Value = -31252012

$ws.Range("H77").Value = 31252886
$ws.Range("I77").Value = 31252886
$ws.Range("K77").Value = 156264430
$ws.Range("M77").Value = -156260062

$ws.Range("H88").Value = 2368.875
$ws.Range("J88").Value = 2345.3
$ws.Range("L88").Value = 2345.3
$ws.Range("N88").Value = -3157.3

$ws.Range("H91").Value = 2368.875
$ws.Range("J91").Value = 2345.3
$ws.Range("L91").Value = 2345.3
$ws.Range("N91").Value = -5153.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2226.0264
$ws.Range("I86").Value = 2422.6086
$ws.Range("J86").Value = 1924.6
$ws.Range("K86").Value = 2422.6086
$ws.Range("L86").Value = 1924.6
$ws.Range("M86").Value = -1299.6086
$ws.Range("N86").Value = -4170.6

$ws.Range("H89").Value = 2226.0264
$ws.Range("I89").Value = 2422.6086
$ws.Range("J89").Value = 1924.6
$ws.Range("K89").Value = 12113.043
$ws.Range("L89").Value = 9623
$ws.Range("M89").Value = -6497.043
$ws.Range("N89").Value = -20855

$ws.Range("H94").Value = 1409.0416
$ws.Range("I94").Value = 1688.0555
$ws.Range("J94").Value = 572
$ws.Range("K94").Value = 1688.0555
$ws.Range("L94").Value = 572
$ws.Range("M94").Value = -1237.0555
$ws.Range("N94").Value = -1474

$ws.Range("H105").Value = 2690
$ws.Range("I105").Value = 2274.182
$ws.Range("K105").Value = 2274.182
$ws.Range("M105").Value = -527.1819999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 110.933334
$ws.Range("I7").Value = 125.5
$ws.Range("J7").Value = 52.666668
$ws.Range("K7").Value = 125.5
$ws.Range("L7").Value = 52.666668
$ws.Range("M7").Value = -12.5
$ws.Range("N7").Value = -278.666668

$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 999
$ws.Range("M22").Value = -649

$ws.Range("H31").Value = 3984.7
$ws.Range("I31").Value = 4983
$ws.Range("K31").Value = 4983
$ws.Range("M31").Value = -4688

$ws.Range("H34").Value = 3984.7
$ws.Range("I34").Value = 4983
$ws.Range("K34").Value = 4983
$ws.Range("M34").Value = -4781

$ws.Range("H62").Value = 2221.6667
$ws.Range("I62").Value = 1839.6
$ws.Range("K62").Value = 1839.6
$ws.Range("M62").Value = -1215.6

$ws.Range("H65").Value = 2221.6667
$ws.Range("I65").Value = 1839.6
$ws.Range("K65").Value = 9198
$ws.Range("M65").Value = -6078

$ws.Range("H134").Value = 10463012
$ws.Range("I134").Value = 11956897
$ws.Range("K134").Value = 35870691
$ws.Range("M134").Value = -35868156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 9510.637000000001
$ws.Range("I26").Value = 240.28572
$ws.Range("J26").Value = 25733.75
$ws.Range("K26").Value = 720.85716
$ws.Range("L26").Value = 77201.25
$ws.Range("M26").Value = -432.85716
$ws.Range("N26").Value = -77777.25

$ws.Range("H81").Value = 289088.78
$ws.Range("I81").Value = 306474.88
$ws.Range("K81").Value = 919424.64
$ws.Range("M81").Value = -918301.64

$ws.Range("H84").Value = 289088.78
$ws.Range("I84").Value = 306474.88
$ws.Range("K84").Value = 2758273.92
$ws.Range("M84").Value = -2752657.92

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H107").Value = 2688.3462
$ws.Range("I107").Value = 2824.1538
$ws.Range("J107").Value = 2552.5386
$ws.Range("K107").Value = 2824.1538
$ws.Range("L107").Value = 2552.5386
$ws.Range("M107").Value = -904.1538
$ws.Range("N107").Value = -6392.5386

$ws.Range("H132").Value = 7357137.5
$ws.Range("I132").Value = 7816333.5
$ws.Range("K132").Value = 23449000.5
$ws.Range("M132").Value = -23446470.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2493.9167
$ws.Range("I40").Value = 2493.9167
$ws.Range("K40").Value = 2493.9167
$ws.Range("M40").Value = -2357.9167

$ws.Range("H122").Value = 5939.7144
$ws.Range("I122").Value = 5600.6665
$ws.Range("K122").Value = 16801.9995
$ws.Range("M122").Value = -14351.9995

$ws.Range("H132").Value = 9263803
$ws.Range("I132").Value = 9619950
$ws.Range("K132").Value = 28859850
$ws.Range("M132").Value = -28857320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 33900
$ws.Range("I43").Value = 33900
$ws.Range("K43").Value = 33900
$ws.Range("M43").Value = -33751

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H81").Value = 1362.7142
$ws.Range("I81").Value = 756.5
$ws.Range("K81").Value = 1513
$ws.Range("M81").Value = -452

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H84").Value = 1362.7142
$ws.Range("I84").Value = 756.5
$ws.Range("K84").Value = 7565
$ws.Range("M84").Value = -2261

$ws.Range("H96").Value = 3218
$ws.Range("I96").Value = 987.375
$ws.Range("K96").Value = 987.375
$ws.Range("M96").Value = 385.625

$ws.Range("H100").Value = 1977.1666
$ws.Range("I100").Value = 2131
$ws.Range("J100").Value = 285
$ws.Range("K100").Value = 4262
$ws.Range("L100").Value = 570
$ws.Range("M100").Value = -3721
$ws.Range("N100").Value = -1652

$ws.Range("H107").Value = 522.55554
$ws.Range("I107").Value = 526.75
$ws.Range("J107").Value = 489
$ws.Range("K107").Value = 1580.25
$ws.Range("L107").Value = 1467
$ws.Range("M107").Value = 339.75
$ws.Range("N107").Value = -5307

$ws.Range("H126").Value = 2977.6
$ws.Range("I126").Value = 2659.8333
$ws.Range("J126").Value = 4248.6665
$ws.Range("K126").Value = 7979.499899999999
$ws.Range("L126").Value = 12745.9995
$ws.Range("M126").Value = -5509.499899999999
$ws.Range("N126").Value = -17685.9995
